# The deck's single addressable theme (ppt/theme/theme1.xml, bound to the
# slide master) switches from the "Integral" palette to the stock
# "Office Theme" palette. dk1/lt1 (black/white) are already correct, so
# only the remaining ten theme-color slots need to change.
#
# PowerPoint's ThemeColorScheme.Colors(index).RGB uses the Windows
# BGR-packed OLE_COLOR encoding (0xBBGGRR), so each target sRGB value is
# byte-reversed below.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

# index -> (slot name, BGR-packed target color)
$targets = @{
    3  = 0x6A5444   # dk2      -> 44546A
    4  = 0xE6E6E7   # lt2      -> E7E6E6
    5  = 0xD59B5B   # accent1  -> 5B9BD5
    6  = 0x317DED   # accent2  -> ED7D31
    7  = 0xA5A5A5   # accent3  -> A5A5A5
    8  = 0x00C0FF   # accent4  -> FFC000
    9  = 0xC47244   # accent5  -> 4472C4
    10 = 0x47AD70   # accent6  -> 70AD47
    11 = 0xC16305   # hlink    -> 0563C1
    12 = 0x724F95   # folHlink -> 954F72
}

foreach ($index in $targets.Keys) {
    $tcs.Colors($index).RGB = $targets[$index]
}
